$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header value changes
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 changes: D2 deleted, B2 and C2 added
$ws.Range("D2").Value = $null
$ws.Range("B2").Value = 28.368239442344922
$ws.Range("C2").Value = 25.532773574620265

# Row 3 changes: B3 deleted, C3 updated
$ws.Range("B3").Value = $null
$ws.Range("C3").Value = 24.181225152760472

# Update selection to reflect new range
$ws.Range("B1:E3").Select()
